$wb = $excel.ActiveWorkbook

# --- Update the dividend data on the "Yearly" sheet ---
# D13 (November dividend) changes from 89.57 to 102.41; every dependent
# formula (G13, D15, G15, and the mirrored totals on "All Time") recalculates
# automatically.
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("D13").Value = 102.41

# Update the selection/view state recorded for the "Yearly" sheet.
$wsYearly.Range("J12").Select()

# --- Update the selection/view state on the "All Time" sheet ---
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Activate()

# Scroll the window so row 31 is at the top of the viewport.
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1

$wsAllTime.Range("I13").Select()
